$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Prediction column (C) values to reflect new picks
$ws.Range("C2").Value = "Marquel Mederos"
$ws.Range("C9").Value = "Waldo Cortes Acosta"

# Move the active selection to C14
$ws.Range("C14").Select()
